$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing A5 timestamp value (tiny floating point correction)
$ws.Range("A5").Value = 45866.2502795949

# Append new row 6 with the latest sensor reading
$ws.Range("A6").Value = 45866.291911875
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat

$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = 31
$ws.Range("D6").Value = 13.67
$ws.Range("E6").Value = 90.88
$ws.Range("F6").Value = 28.13
$ws.Range("G6").Value = 5.8
$ws.Range("H6").Value = "NE"
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "07:00:21"
